# "Playing around with the auto plots a little"
# - Row 12 ("Altitude" plot) is repurposed into the "Des Acceleration" plot.
# - A new row 13 is added for the matching "Acceleration" (actual) plot,
#   mirroring the existing Des-X / X pairs used for Roll/Pitch/Yaw rate above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 scaffolding: numbers / unaffected text first -------------------
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 2

# D13 needs the same "S" label + centred style as D12 - easiest to just
# duplicate the existing formatted cell.
$ws.Range("D12").Copy($ws.Range("D13"))

$ws.Range("E13").Value = "Time [ s ]"

# --- Column F (axis label) : row 12 retargeted, row 13 added together ------
$ws.Range("F12").Value = "Acceleration~[~m/s/s~]"
$ws.Range("F13").Value = "Acceleration~[~m/s/s~]"

$ws.Range("G13").Value = "Vertical"

# --- Column K (channel) : row 12 retargeted, row 13 added together ---------
$ws.Range("K12").Value = "RATE/ADes"
$ws.Range("K13").Value = "RATE/A"

$ws.Range("O13").Value = 1

# --- Column R (override label) : row 12 retargeted, row 13 added together --
$ws.Range("R12").Value = "Des Acceleration"
$ws.Range("R13").Value = "Acceleration"

# --- Column S (override unit) : row 12 retargeted, row 13 added together ---
$ws.Range("S12").Value = "m/s/s"
$ws.Range("S13").Value = "m/s/s"

# Leave the cursor where the author apparently left it when done tinkering.
$ws.Range("S14").Select()
